$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.901.00'
$ws.Range('E2').Value = '  +6.05%  '

$ws.Range('D3').Value = '3.654.87'
$ws.Range('E3').Value = '  +17.75%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '618.87'
$ws.Range('E5').Value = '  +7.36%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.35'
$ws.Range('E6').Value = '  +2.15%  '

$ws.Range('D7').Value = '3.651.36'
$ws.Range('E7').Value = '  +17.68%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('E9').Value = '  +5.22%  '

$ws.Range('E10').Value = '  +8.29%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.67'
$ws.Range('E11').Value = '  +5.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.40'
$ws.Range('E13').Value = '  +11.81%  '

$ws.Range('E14').Value = '  +5.82%  '

$ws.Range('D15').Value = '4.264.84'
$ws.Range('E15').Value = '  +17.82%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.922.26'
$ws.Range('E16').Value = '  +6.09%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.660.04'
$ws.Range('E17').Value = '  +18.01%  '

$ws.Range('E18').Value = '  +1.94%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.56'
$ws.Range('E19').Value = '  +7.42%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '520.87'
$ws.Range('E20').Value = '  +8.42%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.91'
$ws.Range('E21').Value = '  +1.29%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.28'
$ws.Range('E22').Value = '  +18.56%  '

$ws.Range('E23').Value = '  +7.78%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.80'
$ws.Range('E24').Value = '  +6.14%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.53'
$ws.Range('E25').Value = '  +12.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.46'
$ws.Range('E26').Value = '  +6.74%  '

$ws.Range('E27').Value = '  +9.75%  '

$ws.Range('E28').Value = '  -0.06%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.55'
$ws.Range('E29').Value = '  +11.59%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.15'
$ws.Range('E30').Value = '  +3.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.91'
$ws.Range('E31').Value = '  +11.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.62'
$ws.Range('E32').Value = '  +13.05%  '

$ws.Range('E33').Value = '  +17.27%  '

$ws.Range('E34').Value = '  +3.91%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.00%  '

$ws.Range('E36').Value = '  +9.47%  '

$ws.Range('E37').Value = '  +9.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.349'
$ws.Range('E38').Value = '  +11.78%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.19'
$ws.Range('E39').Value = '  +9.53%  '

$ws.Range('E40').Value = '  +7.07%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.67'
$ws.Range('E41').Value = '  +5.33%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '45.81'
$ws.Range('E42').Value = '  -5.08%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.82'
$ws.Range('E43').Value = '  +5.85%  '

$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '426.37'
$ws.Range('E44').Value = '  +13.78%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '3.115.52'
$ws.Range('E45').Value = '  +11.22%  '

$ws.Range('E46').Value = '  +3.58%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0370'
$ws.Range('E47').Value = '  +7.81%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '28.49'
$ws.Range('E48').Value = '  +11.64%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.74'
$ws.Range('E49').Value = '  +4.03%  '

$ws.Range('E50').Value = '  +0.01%  '

$ws.Range('E51').Value = '  +9.99%  '
